# DiscountModel.xlsx edit: add two new discount rows (Thịt bò / Thịt cá)
# and bump the "Thịt heo" row's Approximate Day from 4 to 1.75.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (existing "Thịt heo" row): Approximate Day 4 -> 1.75
$ws.Range("C2").Value = 1.75

# New row 3: a fresh UserID, ProductType "Thịt bò", Approximate Day 4
$ws.Range("A3").Value = "d8298c9e-5d5e-458e-abe0-052340530ce6"
$ws.Range("B3").Value = "Thịt bò"
$ws.Range("C3").Value = 4

# New row 4: same UserID as row 2, ProductType "Thịt cá", Approximate Day 4
$ws.Range("A4").Value = "db08ecfa-c97d-41c3-b306-0b411382c96a"
$ws.Range("B4").Value = "Thịt cá"
$ws.Range("C4").Value = 4

# Column A carries the bold/centered/bordered header-ish style (s="1") on every
# data row -- copy that formatting from the existing A2 cell onto the two new
# UserID cells so they match the rest of the column.
$ws.Range("A2").Copy()
$ws.Range("A3:A4").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Leave the selection where the author's saved file shows it.
$ws.Range("D9").Select() | Out-Null
